# BOT; UPDATE DATA
# Appends one new day's row (2020-05-28, serial 43979) of COVID-19 PCR
# testing data to each of the three data sheets ("all", "kobe", "other"),
# and moves the "active sheet" focus from "kobe" to "other" (the sheet
# the bot visits last in its daily update run).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": simple append after the last row (50 -> 51), formats
# copied from the row above so the new row inherits the same look.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(50).Copy()
$wsAll.Range("A51:H51").PasteSpecial(-4122)  # xlPasteFormats
$wsAll.Range("A51").Value = 43979
$wsAll.Range("B51").Value = 285
$wsAll.Range("C51").Value = 282
$wsAll.Range("D51").Value = 15
$wsAll.Range("E51").Value = 12
$wsAll.Range("F51").Value = 3
$wsAll.Range("G51").Value = 12
$wsAll.Range("H51").Value = 255

# ---------------------------------------------------------------------
# Sheet "kobe": the sheet ends with a footnote row, so the new data
# row must be inserted ABOVE it (pushing the footnote from row 106
# down to row 107) rather than simply appended.
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Rows.Item(105).Copy()
$wsKobe.Rows.Item(106).Insert()
$wsKobe.Range("A106").Value = 43979
$wsKobe.Range("B106").Value = 15
$wsKobe.Range("C106").Value = 3121
$wsKobe.Range("E106").Value = 285
$wsKobe.Range("F106").Value = 12
$wsKobe.Range("G106").Value = 10
$wsKobe.Range("H106").Value = 2
$wsKobe.Range("I106").Value = 12
$wsKobe.Range("J106").Value = 244

# ---------------------------------------------------------------------
# Sheet "other": same footnote-at-the-bottom pattern, new data row
# inserted above it (row 81 -> shifts footnote to row 82).
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(80).Copy()
$wsOther.Rows.Item(81).Insert()
$wsOther.Range("A81").Value = 43979
$wsOther.Range("B81").Value = 0
$wsOther.Range("C81").Value = 14
$wsOther.Range("D81").Value = 3
$wsOther.Range("E81").Value = 2
$wsOther.Range("F81").Value = 1
$wsOther.Range("G81").Value = 0
$wsOther.Range("H81").Value = 11

# ---------------------------------------------------------------------
# Selections: leave each sheet's cursor on the freshly written row,
# matching where the bot's script would have left the cursor after
# writing the day's data.
# ---------------------------------------------------------------------
$wsAll.Range("B51:H51").Select()
$wsKobe.Range("B106:J106").Select()
$wsOther.Range("B85").Select()

# The "other" sheet is the last one touched, so it ends up the active
# (visible) tab when the workbook is saved.
$wsOther.Activate()
